# Refresh market-price-derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# for the rows whose Universalis snapshot changed since the last scheduled run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 1357.3846
$ws.Range("I40").Value = 1289.579
$ws.Range("J40").Value = 1541.4286
$ws.Range("K40").Value = 1289.579
$ws.Range("L40").Value = 1541.4286
$ws.Range("M40").Value = -1114.579
$ws.Range("N40").Value = -1891.4286
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Range("H137").Value = 1821.5264
$ws.Range("J137").Value = 2176.375
$ws.Range("L137").Value = 6529.125
$ws.Range("N137").Value = -11629.125
# Row 138: All-night Crafting / Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3628.6292
$ws.Range("I138").Value = 1423.9714
$ws.Range("J138").Value = 6486.5186
$ws.Range("K138").Value = 4271.914199999999
$ws.Range("L138").Value = 19459.5558
$ws.Range("M138").Value = 868.0858000000007
$ws.Range("N138").Value = -29739.5558

$ws = $wb.Worksheets.Item("ARM")
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 3464.4517
$ws.Range("I61").Value = 3146.8572
$ws.Range("J61").Value = 4131.4
$ws.Range("K61").Value = 3146.8572
$ws.Range("L61").Value = 4131.4
$ws.Range("M61").Value = -2934.8572
$ws.Range("N61").Value = -4555.4
# Row 68: Let Faith Light the Way / Mythrite Bladed Lantern Shield
$ws.Range("H68").Value = 12000
$ws.Range("J68").Value = 12000
$ws.Range("L68").Value = 12000
$ws.Range("N68").Value = -13622
# Row 71: Fifty Shields of Blades (L) / Mythrite Bladed Lantern Shield
$ws.Range("H71").Value = 12000
$ws.Range("J71").Value = 12000
$ws.Range("L71").Value = 36000
$ws.Range("N71").Value = -44112
# Row 82: Belle of the Brawl / Titanium Vambraces of Fending
$ws.Range("H82").Value = 28600
$ws.Range("J82").Value = 28600
$ws.Range("L82").Value = 28600
$ws.Range("N82").Value = -29322
# Row 85: Shouldering the Shut-ins (L) / Titanium Vambraces of Fending
$ws.Range("H85").Value = 28600
$ws.Range("J85").Value = 28600
$ws.Range("L85").Value = 28600
$ws.Range("N85").Value = -31096
# Row 87: Look Before You Leap / Adamantite Leg Guards of Maiming
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
# Row 90: Patience, Young Grasshopper (L) / Adamantite Leg Guards of Maiming
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 3491.818
$ws.Range("I102").Value = 2738.3333
$ws.Range("J102").Value = 4396
$ws.Range("K102").Value = 2738.3333
$ws.Range("L102").Value = 4396
$ws.Range("M102").Value = -1116.3333
$ws.Range("N102").Value = -7640
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 3464.4517
$ws.Range("I136").Value = 3146.8572
$ws.Range("J136").Value = 4131.4
$ws.Range("K136").Value = 9440.5716
$ws.Range("L136").Value = 12394.2
$ws.Range("M136").Value = -6890.571599999999
$ws.Range("N136").Value = -17494.2

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 79.75
$ws.Range("I7").Value = 52
$ws.Range("J7").Value = 99.57143
$ws.Range("K7").Value = 52
$ws.Range("L7").Value = 99.57143
$ws.Range("M7").Value = 61
$ws.Range("N7").Value = -325.57143
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 3153998.5
$ws.Range("I31").Value = 3761317.8
$ws.Range("J31").Value = 269232
$ws.Range("K31").Value = 3761317.8
$ws.Range("L31").Value = 269232
$ws.Range("M31").Value = -3761022.8
$ws.Range("N31").Value = -269822
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 3153998.5
$ws.Range("I34").Value = 3761317.8
$ws.Range("J34").Value = 269232
$ws.Range("K34").Value = 3761317.8
$ws.Range("L34").Value = 269232
$ws.Range("M34").Value = -3761115.8
$ws.Range("N34").Value = -269636
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 3099
$ws.Range("I62").Value = 2952.5
$ws.Range("J62").Value = 3196.6667
$ws.Range("K62").Value = 2952.5
$ws.Range("L62").Value = 3196.6667
$ws.Range("M62").Value = -2328.5
$ws.Range("N62").Value = -4444.6667
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 3099
$ws.Range("I65").Value = 2952.5
$ws.Range("J65").Value = 3196.6667
$ws.Range("K65").Value = 14762.5
$ws.Range("L65").Value = 15983.3335
$ws.Range("M65").Value = -11642.5
$ws.Range("N65").Value = -22223.3335

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap / Maple Syrup
$ws.Range("H5").Value = 1409
$ws.Range("J5").Value = 1908.8889
$ws.Range("L5").Value = 5726.6667
$ws.Range("N5").Value = -5950.6667
# Row 135: Not-so-secret Ingredient / Royal Maple Syrup
$ws.Range("H135").Value = 1409
$ws.Range("J135").Value = 1908.8889
$ws.Range("L135").Value = 17180.0001
$ws.Range("N135").Value = -22250.0001

$ws = $wb.Worksheets.Item("GSM")
# Row 43: Get the Green Stuff / Malachite Earrings
$ws.Range("H43").Value = 15450
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 28900
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 28900
$ws.Range("M43").Value = -1849
$ws.Range("N43").Value = -29202
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 11724922
$ws.Range("I70").Value = 17050800
$ws.Range("J70").Value = 7990
$ws.Range("K70").Value = 17050800
$ws.Range("L70").Value = 7990
$ws.Range("M70").Value = -17050530
$ws.Range("N70").Value = -8530
# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 11724922
$ws.Range("I73").Value = 17050800
$ws.Range("J73").Value = 7990
$ws.Range("K73").Value = 17050800
$ws.Range("L73").Value = 7990
$ws.Range("M73").Value = -17049864
$ws.Range("N73").Value = -9862
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 74673.07
$ws.Range("I132").Value = 2524.8333
$ws.Range("J132").Value = 128784.25
$ws.Range("K132").Value = 7574.499899999999
$ws.Range("L132").Value = 386352.75
$ws.Range("M132").Value = -5044.499899999999
$ws.Range("N132").Value = -391412.75

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 35716456
$ws.Range("I122").Value = 2172
$ws.Range("J122").Value = 71430744
$ws.Range("K122").Value = 6516
$ws.Range("L122").Value = 214292232
$ws.Range("M122").Value = -4066
$ws.Range("N122").Value = -214297132

$ws = $wb.Worksheets.Item("WVR")
# Row 34: He's Got Legs / Velveteen Sarouel
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
# Row 82: Investing in the Future / Hallowed Ramie Doublet of Aiming
$ws.Range("H82").Value = 29800
$ws.Range("J82").Value = 29800
$ws.Range("L82").Value = 29800
$ws.Range("N82").Value = -30566
# Row 85: Maids of Honor (L) / Hallowed Ramie Doublet of Aiming
$ws.Range("H85").Value = 29800
$ws.Range("J85").Value = 29800
$ws.Range("L85").Value = 29800
$ws.Range("N85").Value = -32452
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 1184.2667
$ws.Range("I122").Value = 1184.2667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3552.800099999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1102.800099999999
$ws.Range("N122").ClearContents()
